$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# A new day's price row is inserted at the top (row 2). Every existing data
# row (originally rows 2..186) shifts down by one (to rows 3..187). The
# oldest row (originally row 186) is preserved by re-appending its data as
# the new last row (187). The brand-new top row reuses the previous top
# row's Basic Price / Circular Date / Circular Link (that day's circular
# hasn't been published yet) but gets a new "Date" value one day later.
# ---------------------------------------------------------------------------

$lastRow = 186
$newLastRow = $lastRow + 1

# Snapshot all existing data rows (2..186) in one bulk read.
$data = $ws.Range("A2:F" + $lastRow).Value2

# Remove every hyperlink up front; they will be rebuilt from scratch once
# all the cell text has been shifted into its new location.
$ws.Hyperlinks.Delete()

# Columns A (Date) and E (Circular Date) hold dd-mm-yyyy / dd.mm.yyyy text
# that Excel would otherwise auto-convert to a date serial on assignment;
# prefixing with an apostrophe forces it to stay literal text, exactly as
# it was stored in the source workbook.
function Set-Text($row, $col, $val) {
    $ws.Cells.Item($row, $col).Value2 = "'" + $val
}

# Write the shifted rows: new row r (3..newLastRow) = old row (r-1).
for ($r = $newLastRow; $r -ge 3; $r--) {
    $i = $r - 2
    Set-Text $r 1 $data[$i, 1]
    $ws.Cells.Item($r, 2).Value2 = $data[$i, 2]
    $ws.Cells.Item($r, 3).Value2 = $data[$i, 3]
    $ws.Cells.Item($r, 4).Value2 = $data[$i, 4]
    Set-Text $r 5 $data[$i, 5]
    $ws.Cells.Item($r, 6).Value2 = $data[$i, 6]
}

# New top row: carries forward yesterday's (old row 2's) price/circular
# info, but is labelled with the new date.
Set-Text 2 1 "14-12-2025"
$ws.Cells.Item(2, 2).Value2 = $data[1, 2]
$ws.Cells.Item(2, 3).Value2 = $data[1, 3]
$ws.Cells.Item(2, 4).Value2 = $data[1, 4]
Set-Text 2 5 $data[1, 5]
$ws.Cells.Item(2, 6).Value2 = $data[1, 6]

# Rebuild every hyperlink in column F (2..newLastRow) from the text now
# sitting in each cell (in this sheet the link text always equals its URL).
for ($r = 2; $r -le $newLastRow; $r++) {
    $target = $ws.Cells.Item($r, 6).Value2
    if ($target -ne $null -and $target -ne "") {
        $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $target)
    }
}
